$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Add($null, $summarySheet)
$q3Sheet.Name = "2022-Q3"

# Copy header style (bold, centered, bordered) from an existing quarter sheet
$tplHeader = $wb.Worksheets.Item("2022-Q2").Range("B1:H1")
$tplHeader.Copy($q3Sheet.Range("B1:H1"))

$q3Sheet.Cells.Item(1,2).Value = "基金代码"
$q3Sheet.Cells.Item(1,3).Value = "基金名称"
$q3Sheet.Cells.Item(1,4).Value = "基金规模"
$q3Sheet.Cells.Item(1,5).Value = "股票总仓位"
$q3Sheet.Cells.Item(1,6).Value = "仓位占比"
$q3Sheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3Sheet.Cells.Item(1,8).Value = "仓位排名"

# Template cell for the bold-centered "index" style used in column A
$aStyleTemplate = $wb.Worksheets.Item("2022-Q2").Cells.Item(2,1)

$q3Rows = @(
    "0`t011479`t广发诚享混合A`t36.88`t87.12`t7.82`t2.8840`t7",
    "1`t011130`t广发兴诚混合C`t25.57`t88.20`t6.82`t1.7439`t8",
    "2`t011121`t广发兴诚混合A`t22.16`t88.20`t6.82`t1.5113`t8",
    "3`t519732`t交银定期支付双息平衡混合`t40.09`t66.00`t3.54`t1.4192`t2",
    "4`t161834`t银华鑫锐灵活配置混合（LOF）A`t36.66`t75.53`t3.77`t1.3821`t3",
    "5`t501022`t银华鑫盛灵活配置混合（LOF）A`t41.87`t72.03`t3.27`t1.3691`t3",
    "6`t014423`t汇丰晋信研究精选混合`t37.71`t94.25`t3.52`t1.3274`t10",
    "7`t005535`t泰信竞争优选灵活配置混合`t11.00`t89.95`t7.85`t0.8635`t1",
    "8`t290006`t泰信蓝筹精选混合`t10.65`t90.14`t7.14`t0.7604`t1",
    "9`t013385`t信澳优势价值混合A`t12.44`t84.28`t4.21`t0.5237`t4",
    "10`t010714`t东方红远见价值混合`t15.24`t94.15`t3.13`t0.4770`t9",
    "11`t000480`t东方红新动力灵活配置混合`t12.63`t78.30`t3.04`t0.3840`t8",
    "12`t180010`t银华优质增长混合`t20.22`t82.15`t1.81`t0.3660`t8",
    "13`t011480`t广发诚享混合C`t4.20`t87.12`t7.82`t0.3284`t7",
    "14`t012370`t银华鑫利一年持有期混合`t9.61`t73.11`t3.16`t0.3037`t3",
    "15`t001564`t东方红京东大数据灵活配置混合`t8.84`t73.95`t3.21`t0.2838`t8",
    "16`t169103`t东方红睿轩三年定期开放灵活配置混合`t11.31`t70.03`t2.40`t0.2714`t8",
    "17`t014048`t银华鑫盛灵活配置混合（LOF）C`t8.19`t72.03`t3.27`t0.2678`t3",
    "18`t011405`t银华稳健增长一年持有期混合`t16.19`t79.38`t1.60`t0.2590`t8",
    "19`t014349`t银华鑫锐灵活配置混合（LOF）C`t4.69`t75.53`t3.77`t0.1768`t3",
    "20`t013393`t信澳价值精选混合A`t3.34`t79.98`t3.55`t0.1186`t6",
    "21`t005416`t鹏华尊惠18个月定期开放混合A`t7.53`t37.52`t1.42`t0.1069`t7",
    "22`t015305`t银华鑫峰混合A`t4.56`t54.54`t2.03`t0.0926`t6",
    "23`t920002`t中金精选股票A`t2.95`t82.28`t2.50`t0.0738`t8",
    "24`t013386`t信澳优势价值混合C`t1.26`t84.28`t4.21`t0.0530`t4",
    "25`t673020`t西部利得成长精选灵活配置混合`t1.42`t86.11`t3.38`t0.0480`t7",
    "26`t014321`t德邦周期精选混合A`t0.62`t92.98`t7.42`t0.0460`t3",
    "27`t930602`t国信价值智选混合`t0.49`t76.70`t8.26`t0.0405`t1",
    "28`t006302`t银华行业轮动混合`t1.98`t80.40`t1.65`t0.0327`t8",
    "29`t009667`t鹏华安庆混合A`t2.12`t38.12`t1.40`t0.0297`t5",
    "30`t003165`t鹏华弘嘉灵活配置混合A`t0.82`t91.09`t3.61`t0.0296`t4",
    "31`t015306`t银华鑫峰混合C`t1.25`t54.54`t2.03`t0.0254`t6",
    "32`t009230`t鹏华安和混合A`t1.80`t39.30`t1.40`t0.0252`t6",
    "33`t011572`t鹏华安荣混合A`t1.50`t39.92`t1.37`t0.0206`t8",
    "34`t009668`t鹏华安庆混合C`t1.32`t38.12`t1.40`t0.0185`t5",
    "35`t013394`t信澳价值精选混合C`t0.38`t79.98`t3.55`t0.0135`t6",
    "36`t009231`t鹏华安和混合C`t0.73`t39.30`t1.40`t0.0102`t6",
    "37`t003166`t鹏华弘嘉灵活配置混合C`t0.25`t91.09`t3.61`t0.0090`t4",
    "38`t005519`t银华混改红利灵活配置混合`t0.54`t79.85`t1.63`t0.0088`t8",
    "39`t005417`t鹏华尊惠18个月定期开放混合C`t0.58`t37.52`t1.42`t0.0082`t7",
    "40`t011573`t鹏华安荣混合C`t0.34`t39.92`t1.37`t0.0047`t8",
    "41`t920922`t中金精选股票C`t0.11`t82.28`t2.50`t0.0028`t8",
    "42`t014322`t德邦周期精选混合C`t0.03`t92.98`t7.42`t0.0022`t3"
)

foreach ($r in $q3Rows) {
    $p = $r -split "`t"
    $rowNum = [int]$p[0] + 2
    $aCell = $q3Sheet.Cells.Item($rowNum, 1)
    $aStyleTemplate.Copy($aCell)
    $aCell.Value = [int]$p[0]
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 2) $p[1]
    $q3Sheet.Cells.Item($rowNum, 3).Value = $p[2]
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 4) $p[3]
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 5) $p[4]
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 6) $p[5]
    Set-TextValue $q3Sheet.Cells.Item($rowNum, 7) $p[6]
    $q3Sheet.Cells.Item($rowNum, 8).Value = [int]$p[7]
}

# ---------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a row for 2022-Q3
#    and shift the existing quarters down by one row.
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

$summary.Cells.Item(2,2).Style = "Normal"
$summary.Cells.Item(2,3).Style = "Normal"
$summary.Cells.Item(2,4).Style = "Normal"

$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 43
$summary.Cells.Item(2,4).Value = 17.72

# The freshly inserted row 2 has no "A" cell yet; seed it with the bold
# centered index style (still present on A3, which held the old A2) before
# re-sequencing the whole index column (A) 0..7 across all 8 data rows.
$summary.Cells.Item(3,1).Copy($summary.Cells.Item(2,1))
for ($i = 0; $i -le 7; $i++) {
    $summary.Cells.Item($i + 2, 1).Value = $i
}

Write-Host "done"
